# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 6469
$ws1.Range("F5").Value  = 393
$ws1.Range("F7").Value  = 6
$ws1.Range("F8").Value  = 534
$ws1.Range("F9").Value  = 91
$ws1.Range("F10").Value = 79
$ws1.Range("F13").Value = 378
$ws1.Range("F14").Value = 948
$ws1.Range("F15").Value = 3172
$ws1.Range("F17").Value = 197
$ws1.Range("F18").Value = 1845
$ws1.Range("F19").Value = 24

# Sheet "全部类型" (fourth sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6469
$ws4.Range("F5").Value  = 393
$ws4.Range("F7").Value  = 6
$ws4.Range("F9").Value  = 534
$ws4.Range("F10").Value = 91
$ws4.Range("F11").Value = 79
$ws4.Range("F14").Value = 378
$ws4.Range("F15").Value = 948
$ws4.Range("F16").Value = 3172
$ws4.Range("F18").Value = 197
$ws4.Range("F19").Value = 1845
$ws4.Range("F20").Value = 24
